# Add pdf version of poster - apply layout shifts and text/paragraph tweaks
# to the "two tree-formats" text box on the single poster slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape #29 (id=84, "CustomShape 18") -------------------------------
# Reflow / widen the text box that explains the two tree save formats,
# justify its paragraphs, and rephrase the opening sentence.
$sh29 = $s.Shapes.Item(29)
if ($sh29.Id -ne 84) { throw "Shape 29: expected id 84, got $($sh29.Id)" }
$sh29.Left = 1233.893
$sh29.Width = 522.502

$tr29 = $sh29.TextFrame.TextRange
# Justify every paragraph in this text box (matches <a:pPr algn="just">
# being added across all five paragraphs of this shape).
$tr29.ParagraphFormat.Alignment = 4
# "Prezentowany " (13 chars incl. trailing space) -> "Stworzony "
# leaving "system obsługuje dwa formaty zapisu drzew:" untouched as its
# own run.
$tr29.Characters(1, 13).Text = "Stworzony "

# --- Shape #33 (id=6, "Obraz 5") ----------------------------------------
$sh33 = $s.Shapes.Item(33)
if ($sh33.Id -ne 6) { throw "Shape 33: expected id 6, got $($sh33.Id)" }
$sh33.Top = 843.9027

# --- Shape #36 (id=13, "Obraz 12") ---------------------------------------
$sh36 = $s.Shapes.Item(36)
if ($sh36.Id -ne 13) { throw "Shape 36: expected id 13, got $($sh36.Id)" }
$sh36.Top = 928.7174

# --- Shape #38 (id=54, "CustomShape 18" - "Mankala" label) --------------
$sh38 = $s.Shapes.Item(38)
if ($sh38.Id -ne 54) { throw "Shape 38: expected id 54, got $($sh38.Id)" }
$sh38.Width = 162.7367

# --- Shape #39 (id=55, "CustomShape 18" - "Szachy" label) ---------------
$sh39 = $s.Shapes.Item(39)
if ($sh39.Id -ne 55) { throw "Shape 39: expected id 55, got $($sh39.Id)" }
$sh39.Width = 160.4802
$sh39.Height = 69.8865

# --- Shape #40 (id=15, "Obraz 14") ---------------------------------------
$sh40 = $s.Shapes.Item(40)
if ($sh40.Id -ne 15) { throw "Shape 40: expected id 15, got $($sh40.Id)" }
$sh40.Top = 2027.3906
